$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the hyperlink text that was split across three runs into a
#    single run: "https://medium.com/@vghadigaokar/what" + "-" +
#    "is-stlc-8e19cc51fe80" -> one run with the full URL text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "https://medium.com/@vghadigaokar/what-is-stlc-8e19cc51fe80",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "https://medium.com/@vghadigaokar/what-is-stlc-8e19cc51fe80", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Highlight (yellow) the "As a test lead ..." list paragraph,
#    including the paragraph mark itself, so that the highlight shows
#    up both on the runs and on the paragraph mark's rPr.
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("As a test lead for web based application,")) {
        $p.Range.Font.HighlightColorIndex = 7   # wdYellow
        break
    }
}

# ------------------------------------------------------------------
# 3 & 4. Relocate the "_GoBack" bookmark from inside the "5. Skill
#    Set" paragraph to immediately after the Compatibility Risks run
#    ("... user experience and functionality."), without splitting
#    any text runs. We do this by temporarily inserting a one
#    character marker right after the target text, wrapping a
#    bookmark with that same name around the marker (Bookmarks.Add
#    relocates any existing bookmark of the same name), and then
#    removing the marker character again - this leaves the bookmark
#    collapsed exactly at the desired position.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("This can lead to inconsistencies in the user experience and functionality.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $rng.End

$insertRng = $d.Range($pos, $pos)
$insertRng.InsertAfter([char]1)

$markerRng = $d.Range($pos, $pos + 1)
$d.Bookmarks.Add("_GoBack", $markerRng)

$markerRng2 = $d.Range($pos, $pos + 1)
$markerRng2.Text = ""

# ------------------------------------------------------------------
# 4b. Merge the two runs of the "5. Skill Set ..." paragraph (which
#    used to be split around the old bookmark location) back into a
#    single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "5. Skill Set: Manual testing relies on the expertise and experience of human testers to identify issues, provide feedback, and perform ad hoc testing based on their knowledge of the application.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5. Skill Set: Manual testing relies on the expertise and experience of human testers to identify issues, provide feedback, and perform ad hoc testing based on their knowledge of the application.",
    2) | Out-Null

# ------------------------------------------------------------------
# 5. Merge the "A" + "utomation Testing:" runs into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Automation Testing:", $true, $false, $false, $false, $false, $true, 1, $false,
    "Automation Testing:", 2) | Out-Null
